# Apply cell value updates for the symbol-list refresh (GitHub Actions run).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Numeric-looking values (Price column D) must be forced to Text so Excel
# keeps them as strings instead of coercing them into numbers.
$numericTextChanges = @(
    @{ Ref = "D2"; Value = "241.99" },
    @{ Ref = "D4"; Value = "5.282" },
    @{ Ref = "D5"; Value = "0.05607" },
    @{ Ref = "D6"; Value = "3.378" },
    @{ Ref = "D7"; Value = "6.382" },
    @{ Ref = "D8"; Value = "0.8074" },
    @{ Ref = "D9"; Value = "0.9001" },
    @{ Ref = "D10"; Value = "0.1429" },
    @{ Ref = "D11"; Value = "0.07291" },
    @{ Ref = "D12"; Value = "0.03241" },
    @{ Ref = "D13"; Value = "0.03050" },
    @{ Ref = "D14"; Value = "0.09284" },
    @{ Ref = "D15"; Value = "3.588" },
    @{ Ref = "D16"; Value = "0.001627" },
    @{ Ref = "D17"; Value = "0.04705" },
    @{ Ref = "D18"; Value = "0.0005822" },
    @{ Ref = "D19"; Value = "0.006356" },
    @{ Ref = "D20"; Value = "0.004975" },
    @{ Ref = "D21"; Value = "0.001045" },
    @{ Ref = "D22"; Value = "0.0001503" },
    @{ Ref = "D25"; Value = "2.096" },
    @{ Ref = "D26"; Value = "0.3251" },
    @{ Ref = "D40"; Value = "0.03890" },
    @{ Ref = "D41"; Value = "0.006955" },
    @{ Ref = "D42"; Value = "0.003407" },
    @{ Ref = "D43"; Value = "0.1032" },
    @{ Ref = "D44"; Value = "0.007543" },
    @{ Ref = "D45"; Value = "0.00005939" },
    @{ Ref = "D49"; Value = "0.05826" },
    @{ Ref = "D50"; Value = "0.00002105" },
)

foreach ($change in $numericTextChanges) {
    $cell = $ws.Range($change.Ref)
    $cell.NumberFormat = "@"
    $cell.Value = $change.Value
}

# Plain text values (Coin, Link, Volume columns).
$textChanges = @(
    @{ Ref = "E18"; Value = "17OneONEWorstin24h" },
    @{ Ref = "B42"; Value = "CEJI" },
    @{ Ref = "C42"; Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji" },
    @{ Ref = "E42"; Value = "41CEJICEJI" },
    @{ Ref = "B43"; Value = "BKEXToken" },
    @{ Ref = "C43"; Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk" },
    @{ Ref = "E43"; Value = "42BKEXTokenBKK" },
    @{ Ref = "E47"; Value = "46ACDXExchangeACXT" },
)

foreach ($change in $textChanges) {
    $ws.Range($change.Ref).Value = $change.Value
}
